# Buffer Overflows.docx - "Work on the paper."
#
# The empty paragraph that used to sit right after the "String Manipulation"
# Heading3 (and right before the "Static Code Analysis" Heading2) gets filled
# in with a new paragraph of body text about safe string-manipulation
# libraries (std::string / Boost.StringAlgo). The trailing "_GoBack" bookmark
# that used to mark the end of the "Clang." paragraph moves along with the
# author's cursor to the end of this newly written paragraph.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (currently sitting right
#    after "Clang." near the end of the Static Code Analysis section).
#    Removing it first avoids any duplicate-name collision once the new
#    bookmark is written further down.
# ------------------------------------------------------------------
try {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
} catch {
    # No pre-existing bookmark - nothing to clean up.
}

# ------------------------------------------------------------------
# 2. Locate the empty paragraph that directly follows the
#    "String Manipulation" heading (Heading 3).
# ------------------------------------------------------------------
$headingRange = $d.Content.Duplicate
$headingRange.Find.Execute("String Manipulation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$headingParagraph = $headingRange.Paragraphs(1)
$targetParagraph = $headingParagraph.Next()
$targetRange = $targetParagraph.Range

# ------------------------------------------------------------------
# 3. Replace that empty paragraph's contents with the new text, split
#    across the same four runs as the authored edit, and re-attach the
#    "_GoBack" bookmark at the very end of the paragraph (just before
#    the paragraph mark).
# ------------------------------------------------------------------
$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve">Safe string manipulation libraries are common for C++. Two of the most popular options are the built-in &lt;string&gt; header and the </w:t></w:r>' +
    '<w:r><w:t>Boost String Algorithms Library. The Boost library</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> can be viewed as an ex</w:t></w:r>' +
    '<w:r><w:t>tension of the &lt;string&gt; header, which fills in some gaps for missing functionality. Both libraries provide type safe, flexible and overrun safe classes and functions when used within the spec.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$targetRange.InsertXML($newParagraphXml)
